$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 0
